# Auto-generated: restore recalculated profit figures per scheduled-runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1802.0385
$ws.Range("I15").Value = 1802.0385
$ws.Range("K15").Value = 5406.1155
$ws.Range("M15").Value = -5237.1155
$ws.Range("H41").Value = 9811.909
$ws.Range("J41").Value = 13005.125
$ws.Range("L41").Value = 13005.125
$ws.Range("N41").Value = -13885.125
$ws.Range("H62").Value = 5661.6665
$ws.Range("I62").Value = 6495
$ws.Range("K62").Value = 6495
$ws.Range("M62").Value = -5871
$ws.Range("H65").Value = 5661.6665
$ws.Range("I65").Value = 6495
$ws.Range("K65").Value = 32475
$ws.Range("M65").Value = -29355
$ws.Range("H132").Value = 19451.182
$ws.Range("I132").Value = 4168.619
$ws.Range("K132").Value = 12505.857
$ws.Range("M132").Value = -9975.857
$ws.Range("H138").Value = 2656.4854
$ws.Range("J138").Value = 4018.182
$ws.Range("L138").Value = 12054.546
$ws.Range("N138").Value = -22334.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14494049
$ws.Range("I32").Value = 15385869
$ws.Range("K32").Value = 15385869
$ws.Range("M32").Value = -15385582
$ws.Range("H52").Value = 45000
$ws.Range("J52").Value = 45000
$ws.Range("L52").Value = 45000
$ws.Range("N52").Value = -45636
$ws.Range("H97").Value = 1035.1818
$ws.Range("I97").Value = 1039.2
$ws.Range("J97").Value = 1026.5714
$ws.Range("K97").Value = 1039.2
$ws.Range("L97").Value = 1026.5714
$ws.Range("M97").Value = -543.2
$ws.Range("N97").Value = -2018.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3254
$ws.Range("I86").Value = 2281.8
$ws.Range("K86").Value = 2281.8
$ws.Range("M86").Value = -1158.8
$ws.Range("H89").Value = 3254
$ws.Range("I89").Value = 2281.8
$ws.Range("K89").Value = 11409
$ws.Range("M89").Value = -5793
$ws.Range("H94").Value = 2019.7059
$ws.Range("I94").Value = 2322.36
$ws.Range("J94").Value = 1179
$ws.Range("K94").Value = 2322.36
$ws.Range("L94").Value = 1179
$ws.Range("M94").Value = -1871.36
$ws.Range("N94").Value = -2081
$ws.Range("H99").Value = 21039.912
$ws.Range("I99").Value = 24454
$ws.Range("K99").Value = 24454
$ws.Range("M99").Value = -22956
$ws.Range("H105").Value = 1962
$ws.Range("I105").Value = 1890.6666
$ws.Range("J105").Value = 2176
$ws.Range("K105").Value = 1890.6666
$ws.Range("L105").Value = 2176
$ws.Range("M105").Value = -143.6666
$ws.Range("N105").Value = -5670
$ws.Range("H132").Value = 65450.547
$ws.Range("J132").Value = 65450.547
$ws.Range("L132").Value = 65450.547
$ws.Range("N132").Value = -75570.54699999999
$ws.Range("H134").Value = 1639.6
$ws.Range("I134").Value = 1462
$ws.Range("K134").Value = 4386
$ws.Range("M134").Value = -1851
$ws.Range("H140").Value = 80000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2472.7036
$ws.Range("I132").Value = 2207.1738
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 6621.5214
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -4091.5214
$ws.Range("N132").Value = -17058.5
$ws.Range("H134").Value = 6082.7334
$ws.Range("I134").Value = 4554.923
$ws.Range("K134").Value = 13664.769
$ws.Range("M134").Value = -11129.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 6481911.5
$ws.Range("J29").Value = 266.5
$ws.Range("L29").Value = 799.5
$ws.Range("N29").Value = -1353.5
$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 5000
$ws.Range("K32").Value = 15000
$ws.Range("M32").Value = -14717
$ws.Range("H34").Value = 1716.6666
$ws.Range("I34").Value = 741.6667
$ws.Range("J34").Value = 3666.6667
$ws.Range("K34").Value = 2225.0001
$ws.Range("L34").Value = 11000.0001
$ws.Range("M34").Value = -2141.0001
$ws.Range("N34").Value = -11168.0001
$ws.Range("H46").Value = 100526
$ws.Range("I46").Value = 333566.66
$ws.Range("J46").Value = 651.4286
$ws.Range("K46").Value = 1000699.98
$ws.Range("L46").Value = 1954.2858
$ws.Range("M46").Value = -1000608.98
$ws.Range("N46").Value = -2136.2858
$ws.Range("H88").Value = 4499.857
$ws.Range("J88").Value = 4999.8335
$ws.Range("L88").Value = 14999.5005
$ws.Range("N88").Value = -15855.5005
$ws.Range("H91").Value = 4499.857
$ws.Range("J91").Value = 4999.8335
$ws.Range("L91").Value = 14999.5005
$ws.Range("N91").Value = -17963.5005
$ws.Range("H113").Value = 1215.909
$ws.Range("I113").Value = 1452.6666
$ws.Range("K113").Value = 4357.9998
$ws.Range("M113").Value = -2187.9998
$ws.Range("H122").Value = 1041
$ws.Range("J122").Value = 1066.3334
$ws.Range("L122").Value = 9597.000599999999
$ws.Range("N122").Value = -14497.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 37500
$ws.Range("J82").Value = 37500
$ws.Range("L82").Value = 37500
$ws.Range("N82").Value = -38266
$ws.Range("H85").Value = 37500
$ws.Range("J85").Value = 37500
$ws.Range("L85").Value = 37500
$ws.Range("N85").Value = -40152
$ws.Range("H132").Value = 7799.9
$ws.Range("I132").Value = 8333.223
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 24999.669
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -22469.669
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11292.5
$ws.Range("I40").Value = 14305.625
$ws.Range("J40").Value = 5266.25
$ws.Range("K40").Value = 14305.625
$ws.Range("L40").Value = 5266.25
$ws.Range("M40").Value = -14169.625
$ws.Range("N40").Value = -5538.25
$ws.Range("H46").Value = 1882.7333
$ws.Range("J46").Value = 2766.3333
$ws.Range("L46").Value = 2766.3333
$ws.Range("N46").Value = -3142.3333
$ws.Range("H55").Value = 469.84616
$ws.Range("I55").Value = 735.6667
$ws.Range("J55").Value = 242
$ws.Range("K55").Value = 735.6667
$ws.Range("L55").Value = 242
$ws.Range("M55").Value = -562.6667
$ws.Range("N55").Value = -588
$ws.Range("H68").Value = 1979.7142
$ws.Range("I68").Value = 2038.909
$ws.Range("J68").Value = 1762.6666
$ws.Range("K68").Value = 2038.909
$ws.Range("L68").Value = 1762.6666
$ws.Range("M68").Value = -1289.909
$ws.Range("N68").Value = -3260.6666
$ws.Range("H71").Value = 1979.7142
$ws.Range("I71").Value = 2038.909
$ws.Range("J71").Value = 1762.6666
$ws.Range("K71").Value = 10194.545
$ws.Range("L71").Value = 8813.333000000001
$ws.Range("M71").Value = -6450.545
$ws.Range("N71").Value = -16301.333
$ws.Range("H93").Value = 2953.697
$ws.Range("I93").Value = 1663.1305
$ws.Range("J93").Value = 5922
$ws.Range("K93").Value = 1663.1305
$ws.Range("L93").Value = 5922
$ws.Range("M93").Value = -415.1305
$ws.Range("N93").Value = -8418
$ws.Range("H132").Value = 2449.75
$ws.Range("I132").Value = 2156.7942
$ws.Range("J132").Value = 3445.8
$ws.Range("K132").Value = 6470.382599999999
$ws.Range("L132").Value = 10337.4
$ws.Range("M132").Value = -3940.382599999999
$ws.Range("N132").Value = -15397.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1622.7142
$ws.Range("I100").Value = 1849.6
$ws.Range("J100").Value = 1055.5
$ws.Range("K100").Value = 3699.2
$ws.Range("L100").Value = 2111
$ws.Range("M100").Value = -3158.2
$ws.Range("N100").Value = -3193
$ws.Range("H107").Value = 995.0909
$ws.Range("I107").Value = 945.8333
$ws.Range("J107").Value = 1054.2
$ws.Range("K107").Value = 2837.4999
$ws.Range("L107").Value = 3162.6
$ws.Range("M107").Value = -917.4998999999998
$ws.Range("N107").Value = -7002.6
$ws.Range("H122").Value = 2016.7778
$ws.Range("I122").Value = 1847.6428
$ws.Range("J122").Value = 2608.75
$ws.Range("K122").Value = 5542.928400000001
$ws.Range("L122").Value = 7826.25
$ws.Range("M122").Value = -3092.928400000001
$ws.Range("N122").Value = -12726.25
$ws.Range("H128").Value = 57499.418
$ws.Range("J128").Value = 57499.418
$ws.Range("L128").Value = 57499.418
$ws.Range("N128").Value = -67459.41800000001
$ws.Range("H132").Value = 1157.3226
$ws.Range("I132").Value = 1157.3226
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3471.9678
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -941.9677999999999
$ws.Range("H135").Value = 62597.945
$ws.Range("J135").Value = 62597.945
$ws.Range("L135").Value = 62597.945
$ws.Range("N135").Value = -72737.94500000001
$ws.Range("H136").Value = 1873.25
$ws.Range("J136").Value = 2824
$ws.Range("L136").Value = 8472
$ws.Range("N136").Value = -13572
$ws.Range("N132").ClearContents()
